$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Batumi (Ajara region) gained a II-quarter FDI figure for the preliminary
# 2024 row (row 20) - fill it in and refresh the yearly total in column F.
$ws.Range("C20").Value = 15631.674000000008

$total = $ws.Range("B20").Value2 + $ws.Range("C20").Value2 + $ws.Range("D20").Value2 + $ws.Range("E20").Value2
$ws.Range("F20").Value = $total

$ws.Range("F20").Select()
